$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.155.72"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "1.657.91"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("E4").Value = "  -0.41%  "

$ws.Range("D5").Value = "'218.12"
$ws.Range("E5").Value = "  +0.05%  "

$ws.Range("D6").Value = "'0.5293"
$ws.Range("E6").Value = "  +1.51%  "

$ws.Range("E7").Value = "  -0.32%  "

$ws.Range("D8").Value = "'0.2611"
$ws.Range("E8").Value = "  -1.47%  "

$ws.Range("D9").Value = "'0.06358"
$ws.Range("E9").Value = "  +1.33%  "

$ws.Range("D10").Value = "'20.47"
$ws.Range("E10").Value = "  -1.50%  "

$ws.Range("D11").Value = "'0.07785"
$ws.Range("E11").Value = "  +0.88%  "

$ws.Range("D12").Value = "'4.514"
$ws.Range("E12").Value = "  +1.96%  "

$ws.Range("D13").Value = "1.635.43"
$ws.Range("E13").Value = "  -1.70%  "

$ws.Range("D14").Value = "'0.5493"
$ws.Range("E14").Value = "  +0.95%  "

$ws.Range("D15").Value = "0.0₅8216"
$ws.Range("E15").Value = "  +0.64%  "

$ws.Range("D16").Value = "'65.51"
$ws.Range("E16").Value = "  +1.61%  "

$ws.Range("D17").Value = "26.150.99"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = "  -0.40%  "

$ws.Range("D19").Value = "'4.581"
$ws.Range("E19").Value = "  -1.48%  "

$ws.Range("D20").Value = "'193.08"
$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").Value = "'10.11"
$ws.Range("E21").Value = "  +0.41%  "

$ws.Range("D22").Value = "'6.057"
$ws.Range("E22").Value = "  +0.47%  "

$ws.Range("E23").Value = "  -0.53%  "

$ws.Range("D24").Value = "'141.61"
$ws.Range("E24").Value = "  +1.39%  "

$ws.Range("D25").Value = "'0.1253"
$ws.Range("E25").Value = "  +1.96%  "

$ws.Range("D26").Value = "'7.284"
$ws.Range("E26").Value = "  +1.98%  "

$ws.Range("E27").Value = "  +0.83%  "

$ws.Range("D28").Value = "'1.441"
$ws.Range("E28").Value = "  +1.77%  "

$ws.Range("D29").Value = "'0.05950"
$ws.Range("E29").Value = "  -2.64%  "

$ws.Range("D30").Value = "'1.281"
$ws.Range("E30").Value = "  +0.20%  "

$ws.Range("D31").Value = "'3.527"
$ws.Range("E31").Value = "  -1.20%  "

$ws.Range("D32").Value = "'3.268"
$ws.Range("E32").Value = "  +0.81%  "

$ws.Range("D33").Value = "'1.586"
$ws.Range("E33").Value = "  -2.04%  "

$ws.Range("D34").Value = "'0.9564"
$ws.Range("E34").Value = "  -0.70%  "

$ws.Range("D35").Value = "'2.791"
$ws.Range("E35").Value = "  +0.32%  "

$ws.Range("D36").Value = "'2.413"
$ws.Range("E36").Value = "  -0.46%  "

$ws.Range("D37").Value = "'0.5722"
$ws.Range("E37").Value = "  +0.79%  "

$ws.Range("D38").Value = "'0.01621"
$ws.Range("E38").Value = "  +1.39%  "

$ws.Range("D39").Value = "'5.813"
$ws.Range("E39").Value = "  -2.85%  "

$ws.Range("D40").Value = "'0.8473"
$ws.Range("E40").Value = "  -1.05%  "

$ws.Range("E41").Value = "  -0.27%  "

$ws.Range("D42").Value = "'103.09"
$ws.Range("E42").Value = "  +3.03%  "

$ws.Range("D43").Value = "1.026.02"
$ws.Range("E43").Value = "  +1.08%  "

$ws.Range("D44").Value = "1.802.34"
$ws.Range("E44").Value = "  -0.18%  "

$ws.Range("D45").Value = "'57.40"
$ws.Range("E45").Value = "  +0.63%  "

$ws.Range("E46").Value = "  +0.30%  "

$ws.Range("E47").Value = "  +0.88%  "

$ws.Range("D48").Value = "'0.4294"
$ws.Range("E48").Value = "  +1.76%  "

$ws.Range("D49").Value = "'0.05155"
$ws.Range("E49").Value = "  -0.60%  "

$ws.Range("D50").Value = "'7.827"
$ws.Range("E50").Value = "  -2.29%  "

$ws.Range("D51").Value = "'0.09725"
$ws.Range("E51").Value = "  +0.40%  "
